# Applies the "Corrected text in boxes" edit:
#   - Rectangle 4: "Trial parameters file" -> "Local parameters file"
#                  "PARAMETER_TRIAL -- ASCII" -> "LOCALPARAM_INFO -- ASCII"
#   - Rectangle 9: "Trial parameters file " text re-split across two runs (unchanged text)
#                  "PARAMETER_TRIAL -- ASCII" -> "PARAMETER_TRIAL -- NC"
#   - Oval 23 / Oval 24 ("3" / "4"): drop the now-superfluous trailing endParaRPr
#     (achieved by clearing + retyping the run so no stray end-of-paragraph
#     run-properties element survives).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Rectangle 4 ("Local parameters file" / "LOCALPARAM_INFO -- ASCII") ---
$rect4 = $s.Shapes.Item(1)
$tr4 = $rect4.TextFrame.TextRange

$para1 = $tr4.Paragraphs(1)
$para1.Text = "Local parameters "
[void]$para1.InsertAfter("file")

$para3 = $tr4.Paragraphs(3)
$para3.Text = "-- ASCII"
[void]$para3.InsertBefore("LOCALPARAM_INFO ")

# --- Rectangle 9 ("Trial parameters file " / "PARAMETER_TRIAL -- NC") ---
$rect9 = $s.Shapes.Item(5)
$tr9 = $rect9.TextFrame.TextRange

$para1b = $tr9.Paragraphs(1)
$para1b.Text = "Trial parameters "
[void]$para1b.InsertAfter("file ")

$para3b = $tr9.Paragraphs(3)
$para3b.Text = "PARAMETER_TRIAL -- "
[void]$para3b.InsertAfter("NC")

# --- Oval 23 ("3") / Oval 24 ("4"): strip the trailing endParaRPr ---
$oval23 = $s.Shapes.Item(10)
$tr23 = $oval23.TextFrame.TextRange
$oval23Text = $tr23.Text
$tr23.Delete()
$tr23.Text = $oval23Text

$oval24 = $s.Shapes.Item(11)
$tr24 = $oval24.TextFrame.TextRange
$oval24Text = $tr24.Text
$tr24.Delete()
$tr24.Text = $oval24Text
